# Expand RSD process sets and tables
# - buildings (apartments, attached, detached)
# - space heating/cooling and hot water supply

$wb = $excel.ActiveWorkbook

# --- IND_Sets-Proc: update the remembered selection (no data change there).
$wsInd = $wb.Worksheets.Item("IND_Sets-Proc")
$wsInd.Range("H15:K15").Select()

$ws = $wb.Worksheets.Item("RSD_Sets-Proc")

# --- New data rows 13-26: write B/F/G cells in the exact order needed to reproduce
# --- the shared-string table ordering recorded in the target workbook.
$ws.Cells.Item(13,"G").Value = "Residential - Buildings - Apartments"
$ws.Cells.Item(14,"G").Value = "Residential - Buildings - Attached"
$ws.Cells.Item(15,"G").Value = "Residential - Buildings - Detached"
$ws.Cells.Item(13,"F").Value = "DMD_RSD-BLD-APT"
$ws.Cells.Item(14,"F").Value = "DMD_RSD-BLD-ATT"
$ws.Cells.Item(15,"F").Value = "DMD_RSD-BLD-DET"
$ws.Cells.Item(13,"B").Value = "R*BLD*APT*"
$ws.Cells.Item(14,"B").Value = "R*BLD*ATT*"
$ws.Cells.Item(15,"B").Value = "R*BLD*DET*"
$ws.Cells.Item(16,"B").Value = "R*SW*N*,-R*HP*,-R*HET*"
$ws.Cells.Item(17,"B").Value = "R*SW*HP*N*"
$ws.Cells.Item(18,"B").Value = "R*SW*HET*N*"
$ws.Cells.Item(19,"B").Value = "R*SH*N*,-R*HP*"
$ws.Cells.Item(21,"B").Value = "R*HC*HP*N*"
$ws.Cells.Item(22,"B").Value = "R*SC*N*"
$ws.Cells.Item(23,"B").Value = "R*WH*N*,-R*WH*SOL*"
$ws.Cells.Item(24,"B").Value = "R*WH*SOL*N*"
$ws.Cells.Item(25,"B").Value = "R*SH*X*"
$ws.Cells.Item(26,"B").Value = "R*WH*X*"
$ws.Cells.Item(16,"G").Value = "Residential - Space and Water Heating - New Boilers"
$ws.Cells.Item(17,"G").Value = "Residential - Space and Water Heating - New Heat Pumps"
$ws.Cells.Item(18,"G").Value = "Residential - Space and Water Heating - New District Heating"
$ws.Cells.Item(19,"G").Value = "Residential - Space Heating - New Boilers"
$ws.Cells.Item(20,"G").Value = "Residential - Space Heating - New Heat Pumps"
$ws.Cells.Item(21,"G").Value = "Residential - Space Heating and Cooling - New Heat Pumps"
$ws.Cells.Item(22,"G").Value = "Residential - Space Cooling - New Air Conditioners"
$ws.Cells.Item(23,"G").Value = "Residential - Water Heating - New Boilers"
$ws.Cells.Item(24,"G").Value = "Residential - Water Heating - New Solar"
$ws.Cells.Item(16,"F").Value = "RSD_TECH-SW-N-BOILERS"
$ws.Cells.Item(17,"F").Value = "RSD_TECH-SW-N-HP"
$ws.Cells.Item(18,"F").Value = "RSD_TECH-SW-N-DH"
$ws.Cells.Item(19,"F").Value = "RSD_TECH-SH-N-BOILERS"
$ws.Cells.Item(20,"F").Value = "RSD_TECH-SH-N-HP"
$ws.Cells.Item(21,"F").Value = "RSD_TECH-HC-N-HP"
$ws.Cells.Item(22,"F").Value = "RSD_TECH-SC-N-AC"
$ws.Cells.Item(23,"F").Value = "RSD_TECH-WH-N-BOILERS"
$ws.Cells.Item(24,"F").Value = "RSD_TECH-WH-N-SOLAR"
$ws.Cells.Item(25,"F").Value = "RSD_TECH-SH-E"
$ws.Cells.Item(26,"F").Value = "RSD_TECH-WH-E"
$ws.Cells.Item(26,"G").Value = "Residential - Water Heating - Existing"
$ws.Cells.Item(25,"G").Value = "Residential - Space Heating - Existing"
$ws.Cells.Item(20,"B").Value = "R*SH*HP*N*"

# --- Remaining cells (re-use existing shared strings "DMD"/"AND"/"OR" - order does not
# --- affect the shared string table, write them row by row for clarity).
$ws.Cells.Item(13,"A").Value = "DMD"
$ws.Cells.Item(13,"H").Value = "AND"
$ws.Cells.Item(13,"I").Value = "OR"
$ws.Cells.Item(13,"J").Value = "AND"
$ws.Cells.Item(13,"K").Value = "OR"
$ws.Cells.Item(14,"A").Value = "DMD"
$ws.Cells.Item(14,"H").Value = "AND"
$ws.Cells.Item(14,"I").Value = "OR"
$ws.Cells.Item(14,"J").Value = "AND"
$ws.Cells.Item(14,"K").Value = "OR"
$ws.Cells.Item(15,"A").Value = "DMD"
$ws.Cells.Item(15,"H").Value = "AND"
$ws.Cells.Item(15,"I").Value = "OR"
$ws.Cells.Item(15,"J").Value = "AND"
$ws.Cells.Item(15,"K").Value = "OR"
$ws.Cells.Item(16,"H").Value = "AND"
$ws.Cells.Item(16,"I").Value = "OR"
$ws.Cells.Item(16,"J").Value = "AND"
$ws.Cells.Item(16,"K").Value = "OR"
$ws.Cells.Item(17,"H").Value = "AND"
$ws.Cells.Item(17,"I").Value = "OR"
$ws.Cells.Item(17,"J").Value = "AND"
$ws.Cells.Item(17,"K").Value = "OR"
$ws.Cells.Item(18,"H").Value = "AND"
$ws.Cells.Item(18,"I").Value = "OR"
$ws.Cells.Item(18,"J").Value = "AND"
$ws.Cells.Item(18,"K").Value = "OR"
$ws.Cells.Item(19,"H").Value = "AND"
$ws.Cells.Item(19,"I").Value = "OR"
$ws.Cells.Item(19,"J").Value = "AND"
$ws.Cells.Item(19,"K").Value = "OR"
$ws.Cells.Item(20,"H").Value = "AND"
$ws.Cells.Item(20,"I").Value = "OR"
$ws.Cells.Item(20,"J").Value = "AND"
$ws.Cells.Item(20,"K").Value = "OR"
$ws.Cells.Item(21,"H").Value = "AND"
$ws.Cells.Item(21,"I").Value = "OR"
$ws.Cells.Item(21,"J").Value = "AND"
$ws.Cells.Item(21,"K").Value = "OR"
$ws.Cells.Item(22,"H").Value = "AND"
$ws.Cells.Item(22,"I").Value = "OR"
$ws.Cells.Item(22,"J").Value = "AND"
$ws.Cells.Item(22,"K").Value = "OR"
$ws.Cells.Item(23,"H").Value = "AND"
$ws.Cells.Item(23,"I").Value = "OR"
$ws.Cells.Item(23,"J").Value = "AND"
$ws.Cells.Item(23,"K").Value = "OR"
$ws.Cells.Item(24,"H").Value = "AND"
$ws.Cells.Item(24,"I").Value = "OR"
$ws.Cells.Item(24,"J").Value = "AND"
$ws.Cells.Item(24,"K").Value = "OR"
$ws.Cells.Item(25,"H").Value = "AND"
$ws.Cells.Item(25,"I").Value = "OR"
$ws.Cells.Item(25,"J").Value = "AND"
$ws.Cells.Item(25,"K").Value = "OR"
$ws.Cells.Item(26,"H").Value = "AND"
$ws.Cells.Item(26,"I").Value = "OR"
$ws.Cells.Item(26,"J").Value = "AND"
$ws.Cells.Item(26,"K").Value = "OR"

# --- Column F on RSD_Sets-Proc now holds longer set-name codes (e.g.
# --- "RSD_TECH-SW-N-BOILERS"); widen it to fit, same as the "best fit"
# --- resize Excel performs automatically when the column is selected.
$ws.Columns.Item(6).ColumnWidth = 21.14

# --- Zoom the sheet slightly out and leave the selection/cursor where the
# --- edit session ended, then make this the active sheet/tab.
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("F23").Select()
$ws.Activate()
